$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate "SuplierName" column (H). Excel shifts the
# following "Status" column (previously I) left into H automatically.
$ws.Range("H1").EntireColumn.Delete() | Out-Null

# Update the flight name for row 4 from "Garuda" to "Ramesh"
$ws.Range("G4").Value = "Ramesh"

# Move the active selection to G7, matching the edited workbook's view state
$ws.Range("G7").Select() | Out-Null
